$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 0.2422552835924969
$ws.Range("J4").Value = 0.5268133146755812
$ws.Range("K4").Value = 0.7235849514517669
$ws.Range("L4").Value = 3.022931541449702
